# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K") holds per-game strikeout counts. The values were
# regenerated from the authoritative box-score source ("K" instead of the
# previous "Strike#" derivation), so the corrected counts are written back
# into column G for every affected game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number (sheet row) -> corrected K value
$kValues = @{
    2  = 2
    3  = 2
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 2
    9  = 1
    10 = 2
    11 = 1
    12 = 2
    13 = 2
    14 = 2
    15 = 0
    16 = 1
    17 = 2
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 2
    29 = 0
    30 = 2
    31 = 1
    32 = 1
    33 = 1
    34 = 2
    35 = 3
    36 = 1
    37 = 1
    38 = 0
    39 = 1
    40 = 0
    41 = 2
    42 = 4
    43 = 3
    44 = 1
    45 = 0
    46 = 0
    47 = 0
    48 = 2
    49 = 1
    50 = 1
    51 = 2
    54 = 2
    56 = 1
    57 = 1
}

# Column G is the 7th column ("K" header lives in G1).
$kCol = 7

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, $kCol).Value = $kValues[$row]
}

Write-Host "Updated" $kValues.Keys.Count "K values in column G"
